$wb = $excel.ActiveWorkbook

# --- Metadata sheet: URL + Date updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/environmental-context"
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- Elements sheet: Extension.url fixed value + ValueSet binding URL + column widths ---
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/environmental-context"
$ws2.Range("Z6").Value = "https://2rdoc.pt/fhir/ValueSet/environmental-context"

# Column width adjustments (raw OOXML width = ColumnWidth + 5/6, quantized to the
# nearest 1/6 of a character by the host). Values below are chosen so the stored
# width lands as close as possible to the target from the regenerated export.
$ws2.Columns.Item(1).ColumnWidth = 15.584635416666666
$ws2.Columns.Item(2).ColumnWidth = 15.584635416666666
$ws2.Columns.Item(3).ColumnWidth = 8.959635416666666
$ws2.Columns.Item(3).Hidden = $true
$ws2.Columns.Item(4).ColumnWidth = 6.213541666666667
$ws2.Columns.Item(4).Hidden = $true
$ws2.Columns.Item(5).ColumnWidth = 4.467447916666667
$ws2.Columns.Item(6).ColumnWidth = 3.1197916666666665
$ws2.Columns.Item(7).ColumnWidth = 3.4322916666666665
$ws2.Columns.Item(8).ColumnWidth = 11.854166666666666
$ws2.Columns.Item(9).ColumnWidth = 9.678385416666666
$ws2.Columns.Item(11).ColumnWidth = 13.541666666666666
$ws2.Columns.Item(15).ColumnWidth = 11.428385416666666
$ws2.Columns.Item(20).ColumnWidth = 6.967447916666667
$ws2.Columns.Item(21).ColumnWidth = 12.776041666666666
$ws2.Columns.Item(22).ColumnWidth = 13.084635416666666
$ws2.Columns.Item(23).ColumnWidth = 14.178385416666666
$ws2.Columns.Item(24).ColumnWidth = 13.795572916666666
$ws2.Columns.Item(25).ColumnWidth = 16.248697916666668
$ws2.Columns.Item(26).ColumnWidth = 42.119791666666664
$ws2.Columns.Item(27).ColumnWidth = 4.240885416666667
$ws2.Columns.Item(28).ColumnWidth = 17.147135416666668
$ws2.Columns.Item(29).ColumnWidth = 33.744791666666664
$ws2.Columns.Item(30).ColumnWidth = 12.709635416666666
$ws2.Columns.Item(31).ColumnWidth = 10.486979166666666
$ws2.Columns.Item(31).Hidden = $true
$ws2.Columns.Item(32).ColumnWidth = 14.213541666666666
$ws2.Columns.Item(32).Hidden = $true
$ws2.Columns.Item(33).ColumnWidth = 7.389322916666667
$ws2.Columns.Item(33).Hidden = $true
$ws2.Columns.Item(34).ColumnWidth = 7.697916666666667
$ws2.Columns.Item(37).ColumnWidth = 18.729166666666668
